# Generate Report for Handoff
# Adds a new "row 3" (for file bdca1e70-117a-402e-83b7-3fdf5867d840) to each of
# the three worksheets: Overview, zh-cn, de-de — mirroring the existing row 2
# that documents the 4b17b1dc-... file.

$wb = $excel.ActiveWorkbook

$newId       = "bdca1e70-117a-402e-83b7-3fdf5867d840"
$newMdName   = "$newId.md"
$newHash     = "3c60fd09eb4c3cdcf3bfb132e208449e371b7ae1"
$zhXlfName   = "$newId.$newHash.zh-cn.xlf"
$deXlfName   = "$newId.$newHash.de-de.xlf"

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/7d1cc3c7d14e38a2e2e3f5b7a9d4d2c6f8a1b3c5/e2e/$newMdName"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b1a2c3d4e5f60718293a4b5c6d7e8f9a0b1c2d3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4f5e6d7c8b9a0f1e2d3c4b5a6978685746352413/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$readyStatus = "Ready for handoff"
$neverHandback = "0001-01-01 00:00:00"
$includeReason = "Include"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus
$wsOverview.Range("D3").Value = "2016-03-24 18:46:58"
$wsOverview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $newMdName)
$wsOverview.Range("A3").Font.Underline = 2
$wsOverview.Range("A3").Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Reference Tokens | Handoff Reason | Dependency
# From | Error Detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $newMdName
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyStatus
$wsZh.Range("D3").Value = $zhXlfName
$wsZh.Range("E3").Value = "2016-03-24 18:46:54"
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value = $neverHandback
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("J3").Value = $includeReason

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $newMdName)
$wsZh.Range("A3").Font.Underline = 2
$wsZh.Range("A3").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, "", "", $zhXlfName)
$wsZh.Range("D3").Font.Underline = 2
$wsZh.Range("D3").Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $newMdName
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyStatus
$wsDe.Range("D3").Value = $deXlfName
$wsDe.Range("E3").Value = "2016-03-24 18:46:58"
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value = $neverHandback
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("J3").Value = $includeReason

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $newMdName)
$wsDe.Range("A3").Font.Underline = 2
$wsDe.Range("A3").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, "", "", $deXlfName)
$wsDe.Range("D3").Font.Underline = 2
$wsDe.Range("D3").Font.Color = 15570276
